$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 69 (shifts existing rows 69..157 down to 70..158,
# matching the dimension growing from A1:R157 to A1:R158).
$ws.Rows(69).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(69, 1).Value = 3
$ws.Cells.Item(69, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(69, 3).Value = "Coquimbo"
$ws.Cells.Item(69, 4).Value = 44629
$ws.Cells.Item(69, 5).Value = 5
$ws.Cells.Item(69, 6).Value = 100112030
$ws.Cells.Item(69, 7).Value = "Poroto granado"
$ws.Cells.Item(69, 8).Value = "Sin especificar"
$ws.Cells.Item(69, 9).Value = "Primera"
$ws.Cells.Item(69, 10).Value = 35
$ws.Cells.Item(69, 11).Value = 23000
$ws.Cells.Item(69, 12).Value = 23000
$ws.Cells.Item(69, 13).Value = 23000
$ws.Cells.Item(69, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(69, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(69, 16).Value = 920
$ws.Cells.Item(69, 17).Value = 25
$ws.Cells.Item(69, 18).Value = "Hortaliza"
